# Add a "Save" column (column H) to the s_vals sheet, mirroring the
# existing header style used by the other stat columns (e.g. "sum" in G1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: "Save", styled like the other header cells (G1, etc.)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats - copy style only, keep our value

# Data cells for the new column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
